$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column before column C. This shifts:
#      old C (disk composition formula) -> D
#      old D (spacer)                   -> E
#      old E (notes text)               -> F
#      old H (empty placeholder)        -> I
#      old I (empty placeholder, s=3)   -> J
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).Insert()

# New column C (memory, MB) takes the same visual style as the D column cells
# that used to be in C (style 14) - copy that formatting across C3:C15.
$ws.Range("D3").Copy() | Out-Null
$ws.Range("C3:C15").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# The new C1/C2 header cells should look like their D1/D2 neighbours (style 9).
$ws.Range("D1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 10 + 7/12
$ws.Columns.Item(3).ColumnWidth = 11 + 7/12
$ws.Columns.Item(4).ColumnWidth = 10.25
$ws.Columns.Item(10).ColumnWidth = 11.75

# ---------------------------------------------------------------------------
# 3. Header rows (1 & 2)
# ---------------------------------------------------------------------------
$ws.Range("C1").Value2 = "Mem (MB)"
$ws.Range("D1").Value2 = "Disk (MB)"
$ws.Range("I1").Value2 = "Memory (kb)"
$ws.Range("J1").Value2 = "Memory (MB)"

# ---------------------------------------------------------------------------
# 4. Raw memory values (kb) in column I and formulas (MB) in column J,
#    plus the new memory formula in column C for each data row.
# ---------------------------------------------------------------------------
$memKb = @{
  3  = 89954.4
  4  = 89235.2
  5  = 95281.6
  6  = 86853.6
  7  = 90644.8
  8  = 83992.8
  10 = 221820.8
  11 = 195058.4
  12 = 222985.6
  13 = 197671.2
  14 = 199677.6
  15 = 196940
}

foreach ($row in $memKb.Keys) {
    $ws.Range("I$row").Value2 = $memKb[$row]
    $ws.Range("J$row").Formula = "=I$row/1024"
    $ws.Range("C$row").Formula = "=J$row"
}

# Give the new column J the same number style as column C (style index 2,
# i.e. number format "0.0"). Row 9 is a spacer row with no memory data, so
# leave it untouched (cleared below).
$ws.Range("J1:J8").NumberFormat = "0.0"
$ws.Range("J10:J15").NumberFormat = "0.0"

# ---------------------------------------------------------------------------
# 5. These rows carry no memory data at all - the leftover placeholder cells
#    that used to live in the old I/H columns (now shifted into J) must be
#    removed entirely rather than merely cleared of their value.
# ---------------------------------------------------------------------------
$ws.Range("J9").Clear() | Out-Null
$ws.Range("J16:J29").Clear() | Out-Null

# ---------------------------------------------------------------------------
# 6. Selection moves to F11.
# ---------------------------------------------------------------------------
$ws.Range("F11").Select() | Out-Null
